$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# Fix E28: was stored as text "544028", should be a real number
$ws.Cells.Item(28, 5).Value = 544028

# Append new row 29 with the new screener hit
$ws.Cells.Item(29, 1).Value = "21/06/2024 09:44:57"
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(29, 3).Value = "TATATECH"
$ws.Cells.Item(29, 4).Value = "Tata Technologies Ltd"

# bsecode stays text here (matches the source data quirk), so force it
# with a leading apostrophe and strip the text-number formatting that
# would otherwise get auto-applied.
$ws.Cells.Item(29, 5).Value = "'544028"
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(29, 6).Value = -2.02
$ws.Cells.Item(29, 7).Value = 989.95
$ws.Cells.Item(29, 8).Value = 4168174
